# edit.ps1 — Word COM-interop script reproducing the commit's XML diff:
#   1. Update the rendered "Date: ..." line to the new timestamp.
#   2. Give the document an explicit section (w:sectPr) with the
#      standard A4 / 2.5x2x2.5cm page setup that the regenerated
#      output now carries (page size, margins, column spacing).

$d = $word.ActiveDocument

# --- 1. Replace the stale rendered timestamp with the regenerated one ---
$oldDate = "Date: 2025-07-24 16:12:41.596068 +0200 CEST m=+0.021787834"
$newDate = "Date: 2025-10-20 08:34:24.131949 +0200 CEST m=+0.002295126"

$d.Content.Find.Execute($oldDate, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newDate, 2)

# --- 2. Materialize the document's (previously implicit) section ---
# properties so the last section carries an explicit <w:sectPr> with
# A4 paper, the template's margins, and the default column spacing.
$ps = $d.PageSetup

$ps.PageWidth       = 595.3    # 11906 twips
$ps.PageHeight      = 841.9    # 16838 twips
$ps.TopMargin       = 70.85    # 1417 twips
$ps.RightMargin     = 70.85    # 1417 twips
$ps.BottomMargin    = 56.7     # 1134 twips
$ps.LeftMargin      = 70.85    # 1417 twips
$ps.HeaderDistance  = 35.4     # 708 twips
$ps.FooterDistance  = 35.4     # 708 twips
$ps.Gutter          = 0

$ps.TextColumns.Spacing = 35.4 # 708 twips
